# Generate Report for Handoff
#
# The "b.md" source file has been handed off again (new target files were
# generated for zh-cn and de-de), so:
#   - Overview!B3/C3 status moves from "Handed back: in sync with en-US"
#     to "Ready for handoff"
#   - zh-cn!B3 status moves to "Ready for handoff" too, and the "Latest
#     Handoff File"/"Latest Handoff Datetime" columns (C3/D3) are updated
#     to reflect the new handoff package
#   - de-de!B3/C3/D3 get the equivalent update

$wb = $excel.ActiveWorkbook

$status = "Ready for handoff"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $status
$overview.Range("C3").Value = $status

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = $status
$zhcn.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("D3").Value = "2016-02-18 09:46:49"

# Keep the hyperlink display text on C3 (Latest Handoff File) in sync with
# the new file name - the link target itself is unchanged.
$i = 0
foreach ($hl in $zhcn.Hyperlinks) {
    $i = $i + 1
    if ($i -eq 6) {
        $hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
    }
}

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = $status
$dede.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("D3").Value = "2016-02-18 09:47:00"

$i = 0
foreach ($hl in $dede.Hyperlinks) {
    $i = $i + 1
    if ($i -eq 6) {
        $hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
    }
}
